$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells that would otherwise be auto-converted to numbers
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply updated values
$ws.Range("D2").Value = "41.234.48"
$ws.Range("E2").Value = "  -1.01%  "
$ws.Range("D3").Value = "2.440.09"
$ws.Range("E3").Value = "  -1.36%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "316.85"
$ws.Range("E5").Value = "  -0.58%  "
$ws.Range("D6").Value = "89.31"
$ws.Range("E6").Value = "  -3.88%  "
$ws.Range("E7").Value = "  -2.09%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("D9").Value = "0.497"
$ws.Range("E9").Value = "  -3.93%  "
$ws.Range("D10").Value = "32.21"
$ws.Range("E10").Value = "  -2.47%  "
$ws.Range("D11").Value = "0.0834"
$ws.Range("E11").Value = "  -4.01%  "
$ws.Range("E12").Value = "  -2.73%  "
$ws.Range("D13").Value = "2.809.73"
$ws.Range("E13").Value = "  -1.58%  "
$ws.Range("D14").Value = "6.73"
$ws.Range("E14").Value = "  -2.66%  "
$ws.Range("D15").Value = "15.62"
$ws.Range("E15").Value = "  -0.29%  "
$ws.Range("D16").Value = "2.449.62"
$ws.Range("E16").Value = "  -0.95%  "
$ws.Range("D17").Value = "0.774"
$ws.Range("E17").Value = "  -1.99%  "
$ws.Range("D18").Value = "41.151.97"
$ws.Range("E18").Value = "  -1.12%  "
$ws.Range("E19").Value = "  -3.55%  "
$ws.Range("E20").Value = "  -3.74%  "
$ws.Range("D21").Value = "72.28"
$ws.Range("E21").Value = "  +1.32%  "
$ws.Range("D22").Value = "11.02"
$ws.Range("E22").Value = "  -3.99%  "
$ws.Range("D23").Value = "235.58"
$ws.Range("E23").Value = "  -2.37%  "
$ws.Range("E24").Value = "  -2.36%  "
$ws.Range("E25").Value = "  +0.03%  "
$ws.Range("E26").Value = "  -2.55%  "
$ws.Range("D27").Value = "24.10"
$ws.Range("E27").Value = "  -3.02%  "
$ws.Range("E28").Value = "  -3.02%  "
$ws.Range("D29").Value = "9.56"
$ws.Range("E29").Value = "  -3.49%  "
$ws.Range("D30").Value = "34.92"
$ws.Range("E30").Value = "  -4.75%  "
$ws.Range("D31").Value = "156.45"
$ws.Range("E31").Value = "  -0.48%  "
$ws.Range("E32").Value = "  +0.07%  "
$ws.Range("E33").Value = "  -4.82%  "
$ws.Range("E34").Value = "  -2.24%  "
$ws.Range("D35").Value = "0.0746"
$ws.Range("E35").Value = "  -3.41%  "
$ws.Range("E36").Value = "  +0.50%  "
$ws.Range("D37").Value = "16.70"
$ws.Range("E37").Value = "  -4.32%  "
$ws.Range("E38").Value = "  -0.73%  "
$ws.Range("E39").Value = "  -3.37%  "
$ws.Range("D40").Value = "0.101"
$ws.Range("E40").Value = "  -2.62%  "
$ws.Range("D41").Value = "3.86"
$ws.Range("E41").Value = "  -3.97%  "
$ws.Range("E42").Value = "  -7.03%  "
$ws.Range("D43").Value = "1.991.76"
$ws.Range("E43").Value = "  +0.33%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "18.46"
$ws.Range("E44").Value = "  -4.47%  "
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").Value = "0.0276"
$ws.Range("E45").Value = "  -3.56%  "
$ws.Range("D46").Value = "2.89"
$ws.Range("E46").Value = "  -4.68%  "
$ws.Range("D47").Value = "9.43"
$ws.Range("E47").Value = "  +1.99%  "
$ws.Range("D48").Value = "2.669.85"
$ws.Range("E48").Value = "  -1.58%  "
$ws.Range("D49").Value = "95.36"
$ws.Range("D50").Value = "73.26"
$ws.Range("E50").Value = "  -0.81%  "
$ws.Range("D51").Value = "51.72"
$ws.Range("E51").Value = "  -1.84%  "

Write-Host "Applied cryptos update"